$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.669.83'
$ws.Range("E2").Value = '  +3.82%  '

$ws.Range("D3").Value = '1.911.08'
$ws.Range("E3").Value = '  +1.91%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.73%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.45'
$ws.Range("E5").Value = '  +0.79%  '

$ws.Range("E6").Value = '  -0.76%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5186'
$ws.Range("E7").Value = '  +1.67%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3973'
$ws.Range("E8").Value = '  +1.19%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08492'
$ws.Range("E9").Value = '  +1.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.91'
$ws.Range("E10").Value = '  +2.74%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.122'
$ws.Range("E11").Value = '  +1.00%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.291'
$ws.Range("E12").Value = '  +0.33%  '

$ws.Range("D13").Value = '1.907.15'
$ws.Range("E13").Value = '  +1.41%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.88'
$ws.Range("E14").Value = '  +2.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.346'
$ws.Range("E15").Value = '  +1.11%  '

$ws.Range("E16").Value = '  -0.70%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.94'
$ws.Range("E17").Value = '  +2.75%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001115'
$ws.Range("E18").Value = '  +0.74%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06754'
$ws.Range("E19").Value = '  +0.30%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.96'
$ws.Range("E20").Value = '  +1.30%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  -0.76%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.019'
$ws.Range("E22").Value = '  +0.92%  '

$ws.Range("D23").Value = '29.675.60'
$ws.Range("E23").Value = '  +3.74%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.19'
$ws.Range("E24").Value = '  +0.50%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.208'
$ws.Range("E25").Value = '  -1.79%  '

$ws.Range("D26").Value = '2.124.71'
$ws.Range("E26").Value = '  +1.66%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '159.37'
$ws.Range("E27").Value = '  -1.53%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.99'
$ws.Range("E28").Value = '  +1.12%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.464'
$ws.Range("E29").Value = '  +4.13%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.54'
$ws.Range("E30").Value = '  +1.60%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.084'
$ws.Range("E31").Value = '  +2.73%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1059'
$ws.Range("E32").Value = '  +0.27%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.181'
$ws.Range("E33").Value = '  +6.50%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.677'
$ws.Range("E34").Value = '  +1.81%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02499'
$ws.Range("E35").Value = '  +1.72%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06623'
$ws.Range("E36").Value = '  +1.52%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.146'
$ws.Range("E37").Value = '  +2.82%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2211'
$ws.Range("E38").Value = '  +1.13%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.243'
$ws.Range("E39").Value = '  +4.03%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.186'
$ws.Range("E40").Value = '  +2.41%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6534'
$ws.Range("E41").Value = '  +1.06%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.243'
$ws.Range("E42").Value = '  -1.89%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.44'
$ws.Range("E43").Value = '  +2.32%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6139'
$ws.Range("E44").Value = '  +1.20%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.21'
$ws.Range("E45").Value = '  +1.47%  '

$ws.Range("E46").Value = '  -0.12%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.065'
$ws.Range("E47").Value = '  +1.78%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.238'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.58'
$ws.Range("E49").Value = '  +1.76%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.189'
$ws.Range("E50").Value = '  -0.59%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.44'
$ws.Range("E51").Value = '  +1.93%  '
